$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) -- refreshed "想去人数" (want-to-go count) values in column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 2872
$ws1.Range("F7").Value = 1687
$ws1.Range("F8").Value = 1888
$ws1.Range("F10").Value = 285
$ws1.Range("F11").Value = 770
$ws1.Range("F12").Value = 915
$ws1.Range("F14").Value = 387
$ws1.Range("F15").Value = 1122
$ws1.Range("F17").Value = 53
$ws1.Range("F19").Value = 6837
$ws1.Range("F20").Value = 256
$ws1.Range("F21").Value = 1647
$ws1.Range("F22").Value = 168
$ws1.Range("F25").Value = 317
$ws1.Range("F26").Value = 274
$ws1.Range("F28").Value = 1111
$ws1.Range("F29").Value = 917
$ws1.Range("F30").Value = 62
$ws1.Range("F31").Value = 100
$ws1.Range("F33").Value = 799
$ws1.Range("F34").Value = 1926
$ws1.Range("F35").Value = 163
$ws1.Range("F37").Value = 232
$ws1.Range("F39").Value = 145
$ws1.Range("F40").Value = 233
$ws1.Range("F42").Value = 171

# Sheet "演出" (Performances) -- refreshed "想去人数" (want-to-go count) values in column F
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 16

# Sheet "全部类型" (All types, combined view) -- mirrors the same refreshed values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 16
$ws4.Range("F9").Value = 2872
$ws4.Range("F10").Value = 1687
$ws4.Range("F11").Value = 1888
$ws4.Range("F13").Value = 285
$ws4.Range("F14").Value = 770
$ws4.Range("F16").Value = 915
$ws4.Range("F18").Value = 387
$ws4.Range("F19").Value = 1122
$ws4.Range("F20").Value = 53
$ws4.Range("F22").Value = 6837
$ws4.Range("F23").Value = 256
$ws4.Range("F24").Value = 1647
$ws4.Range("F26").Value = 168
$ws4.Range("F29").Value = 317
$ws4.Range("F30").Value = 274
$ws4.Range("F32").Value = 1111
$ws4.Range("F33").Value = 917
$ws4.Range("F34").Value = 62
$ws4.Range("F35").Value = 100
$ws4.Range("F37").Value = 799
$ws4.Range("F38").Value = 1926
$ws4.Range("F39").Value = 163
$ws4.Range("F41").Value = 232
$ws4.Range("F43").Value = 145
$ws4.Range("F44").Value = 233
$ws4.Range("F49").Value = 171
